# Append 7 new occurrence rows (23-29) for "Tretåig hackspett" (Picoides tridactylus)
# to the Artfynd sheet, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold date-like text (must stay TEXT, not auto-converted to a date serial).
$dateRangeY = $ws.Range("Y23:Y29")
$dateRangeAA = $ws.Range("AA23:AA29")
$dateRangeY.NumberFormat = "@"
$dateRangeAA.NumberFormat = "@"

$rows = @(
    @{
        Row = 23
        A = 131064798; B = 57884; D = "NT"; E = 100109
        F = "Tretåig hackspett"; G = "Picoides tridactylus"; H = "(Linnaeus, 1758)"
        P = "Lars-Olssved, Jmt"; Q = 448209; R = 7037284; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Offerdal"
        Y = "2026-02-07"; AA = "2026-02-07"
        AC = "Ringhack äldre"
        AD = $false; AE = $false; AG = $false
        AW = "Benny Öwre"; AX = "Benny Öwre"
    },
    @{
        Row = 24
        A = 131064797; B = 57884; D = "NT"; E = 100109
        F = "Tretåig hackspett"; G = "Picoides tridactylus"; H = "(Linnaeus, 1758)"
        P = "Lars-Olssved, Jmt"; Q = 448211; R = 7037286; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Offerdal"
        Y = "2026-02-07"; AA = "2026-02-07"
        AC = "Ringhack äldre"
        AD = $false; AE = $false; AG = $false
        AW = "Benny Öwre"; AX = "Benny Öwre"
    },
    @{
        Row = 25
        A = 131064799; B = 57884; D = "NT"; E = 100109
        F = "Tretåig hackspett"; G = "Picoides tridactylus"; H = "(Linnaeus, 1758)"
        P = "Lars-Olssved, Jmt"; Q = 448242; R = 7037242; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Offerdal"
        Y = "2026-02-07"; AA = "2026-02-07"
        AC = "Ringhack färska och äldre"
        AD = $false; AE = $false; AG = $false
        AW = "Benny Öwre"; AX = "Benny Öwre"
    },
    @{
        Row = 26
        A = 131064804; B = 57884; D = "NT"; E = 100109
        F = "Tretåig hackspett"; G = "Picoides tridactylus"; H = "(Linnaeus, 1758)"
        P = "Lars-Olssved, Jmt"; Q = 448308; R = 7037158; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Offerdal"
        Y = "2026-02-07"; AA = "2026-02-07"
        AC = "Ringhack äldre"
        AD = $false; AE = $false; AG = $false
        AW = "Benny Öwre"; AX = "Benny Öwre"
    },
    @{
        Row = 27
        A = 131064796; B = 57884; D = "NT"; E = 100109
        F = "Tretåig hackspett"; G = "Picoides tridactylus"; H = "(Linnaeus, 1758)"
        M = "gammalt bo"
        P = "Lars-Olssved, Jmt"; Q = 448214; R = 7037299; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Offerdal"
        Y = "2026-02-07"; AA = "2026-02-07"
        AC = "Troligt gammalt bohål ca 2,3m upp i granhögstubbe"
        AD = $false; AE = $true; AG = $false
        AW = "Benny Öwre"; AX = "Benny Öwre"
    },
    @{
        Row = 28
        A = 131064802; B = 57884; D = "NT"; E = 100109
        F = "Tretåig hackspett"; G = "Picoides tridactylus"; H = "(Linnaeus, 1758)"
        P = "Lars-Olssved, Jmt"; Q = 448230; R = 7037239; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Offerdal"
        Y = "2026-02-07"; AA = "2026-02-07"
        AC = "Ringhack färska"
        AD = $false; AE = $false; AG = $false
        AW = "Benny Öwre"; AX = "Benny Öwre"
    },
    @{
        Row = 29
        A = 131064800; B = 57884; D = "NT"; E = 100109
        F = "Tretåig hackspett"; G = "Picoides tridactylus"; H = "(Linnaeus, 1758)"
        P = "Lars-Olssved, Jmt"; Q = 448227; R = 7037255; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Offerdal"
        Y = "2026-02-07"; AA = "2026-02-07"
        AC = "Ringhack"
        AD = $false; AE = $false; AG = $false
        AW = "Benny Öwre"; AX = "Benny Öwre"
    }
)

$colOrder = @("A","B","D","E","F","G","H","M","P","Q","R","S","T","U","V","W","Y","AA","AC","AD","AE","AG","AW","AX")

foreach ($rowData in $rows) {
    $r = $rowData.Row
    foreach ($col in $colOrder) {
        if ($rowData.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $rowData[$col]
        }
    }
}

# Restore default (General) formatting for the date-text columns we touched.
$dateRangeY.Style = "Normal"
$dateRangeAA.Style = "Normal"
